# Apply updated pl_mw values for Case_4_140 (380 kV case) to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 3.089818895106305
$ws.Cells.Item(2, 3).Value = 0.6710527065372389
$ws.Cells.Item(2, 4).Value = 0.07649502283095444
$ws.Cells.Item(2, 5).Value = 0.02556318955146963
$ws.Cells.Item(2, 7).Value = 0.002594716094228344
$ws.Cells.Item(2, 9).Value = 2.432396063440251
$ws.Cells.Item(2, 12).Value = 0.2347630007426815
$ws.Cells.Item(2, 13).Value = 0.5549608233648229
$ws.Cells.Item(2, 14).Value = 2.806699538054573
$ws.Cells.Item(3, 2).Value = 2.955433672008326
$ws.Cells.Item(3, 3).Value = 0.6127908413940872
$ws.Cells.Item(3, 4).Value = 0.06966646358756634
$ws.Cells.Item(3, 5).Value = 0.02343348597613826
$ws.Cells.Item(3, 7).Value = 0.002602125036305126
$ws.Cells.Item(3, 9).Value = 2.400344711970476
$ws.Cells.Item(3, 12).Value = 0.2318836519922129
$ws.Cells.Item(3, 13).Value = 0.5356864127814234
$ws.Cells.Item(3, 14).Value = 2.80358162422533
$ws.Cells.Item(4, 2).Value = 2.875088164966996
$ws.Cells.Item(4, 3).Value = 0.5774485357267167
$ws.Cells.Item(4, 4).Value = 0.06552182120543648
$ws.Cells.Item(4, 5).Value = 0.02211642268667724
$ws.Cells.Item(4, 7).Value = 0.002606907520218936
$ws.Cells.Item(4, 9).Value = 2.381854989151776
$ws.Cells.Item(4, 12).Value = 0.2302414597080471
$ws.Cells.Item(4, 13).Value = 0.5242271041360453
$ws.Cells.Item(4, 14).Value = 2.802363001105647
$ws.Cells.Item(5, 2).Value = 2.842886854036522
$ws.Cells.Item(5, 3).Value = 0.563152012576893
$ws.Cells.Item(5, 4).Value = 0.06384465682764073
$ws.Cells.Item(5, 5).Value = 0.0215772361559452
$ws.Cells.Item(5, 7).Value = 0.002608915333195599
$ws.Cells.Item(5, 9).Value = 2.374617358959426
$ws.Cells.Item(5, 12).Value = 0.2296037657318877
$ws.Cells.Item(5, 13).Value = 0.5196510019958609
$ws.Cells.Item(5, 14).Value = 2.802040132278719
$ws.Cells.Item(6, 2).Value = 2.837572332054719
$ws.Cells.Item(6, 3).Value = 0.5607843991617187
$ws.Cells.Item(6, 4).Value = 0.06356686994753602
$ws.Cells.Item(6, 5).Value = 0.02148755227292298
$ws.Cells.Item(6, 7).Value = 0.002609252293864517
$ws.Cells.Item(6, 9).Value = 2.373433437175166
$ws.Cells.Item(6, 12).Value = 0.2294997774063887
$ws.Cells.Item(6, 13).Value = 0.5188967802802296
$ws.Cells.Item(6, 14).Value = 2.801996977233273
$ws.Cells.Item(7, 2).Value = 2.874651706248869
$ws.Cells.Item(7, 3).Value = 0.5772553029716505
$ws.Cells.Item(7, 4).Value = 0.06549915493656044
$ws.Cells.Item(7, 5).Value = 0.02210916117737405
$ws.Cells.Item(7, 7).Value = 0.002606934359530744
$ws.Cells.Item(7, 9).Value = 2.381756179856851
$ws.Cells.Item(7, 12).Value = 0.2302327320866127
$ws.Cells.Item(7, 13).Value = 0.5241650109144587
$ws.Cells.Item(7, 14).Value = 2.802357944969557
$ws.Cells.Item(8, 2).Value = 3.043029956038652
$ws.Cells.Item(8, 3).Value = 0.6508730548909512
$ws.Cells.Item(8, 4).Value = 0.07413037079147955
$ws.Cells.Item(8, 5).Value = 0.02483075834600612
$ws.Cells.Item(8, 7).Value = 0.002597222402242139
$ws.Cells.Item(8, 9).Value = 2.421096475321463
$ws.Cells.Item(8, 12).Value = 0.2337440381652556
$ws.Cells.Item(8, 13).Value = 0.5482366874840849
$ws.Cells.Item(8, 14).Value = 2.805479321366604
$ws.Cells.Item(9, 2).Value = 3.390664095720865
$ws.Cells.Item(9, 3).Value = 0.7987787367351302
$ws.Cells.Item(9, 4).Value = 0.09145203594093232
$ws.Cells.Item(9, 5).Value = 0.03009804859263809
$ws.Cells.Item(9, 7).Value = 0.002580018360902034
$ws.Cells.Item(9, 9).Value = 2.50778819831919
$ws.Cells.Item(9, 12).Value = 0.2416329225504512
$ws.Cells.Item(9, 13).Value = 0.5984540729481154
$ws.Cells.Item(9, 14).Value = 2.81717877627959
$ws.Cells.Item(10, 2).Value = 3.65708244291443
$ws.Cells.Item(10, 3).Value = 0.9097884623059826
$ws.Cells.Item(10, 4).Value = 0.1044405453853159
$ws.Cells.Item(10, 5).Value = 0.03393263077170516
$ws.Cells.Item(10, 7).Value = 0.002568486096492636
$ws.Cells.Item(10, 9).Value = 2.577454640886998
$ws.Cells.Item(10, 12).Value = 0.2480492789715498
$ws.Cells.Item(10, 13).Value = 0.6372399652406386
$ws.Cells.Item(10, 14).Value = 2.829258101989581
$ws.Cells.Item(11, 2).Value = 3.78075948295924
$ws.Cells.Item(11, 3).Value = 0.9608398105238507
$ws.Cells.Item(11, 4).Value = 0.1104109945121365
$ws.Cells.Item(11, 5).Value = 0.03567102171637515
$ws.Cells.Item(11, 7).Value = 0.002563477079565288
$ws.Cells.Item(11, 9).Value = 2.610479695449015
$ws.Cells.Item(11, 12).Value = 0.251104989452358
$ws.Cells.Item(11, 13).Value = 0.655307692988174
$ws.Cells.Item(11, 14).Value = 2.835528160818541
$ws.Cells.Item(12, 2).Value = 3.827956188153337
$ws.Cells.Item(12, 3).Value = 0.9802542791995847
$ws.Cells.Item(12, 4).Value = 0.1126811138785939
$ws.Cells.Item(12, 5).Value = 0.03632857184796023
$ws.Cells.Item(12, 7).Value = 0.002561614138825828
$ws.Cells.Item(12, 9).Value = 2.623179895880014
$ws.Cells.Item(12, 12).Value = 0.25228194311876
$ws.Cells.Item(12, 13).Value = 0.6622113599871398
$ws.Cells.Item(12, 14).Value = 2.838015419514079
$ws.Cells.Item(13, 2).Value = 3.817775296429602
$ws.Cells.Item(13, 3).Value = 0.9760693116345465
$ws.Cells.Item(13, 4).Value = 0.1121917866700954
$ws.Cells.Item(13, 5).Value = 0.03618698768979911
$ws.Cells.Item(13, 7).Value = 0.002562013853931838
$ws.Cells.Item(13, 9).Value = 2.620435995910043
$ws.Cells.Item(13, 12).Value = 0.2520275812090063
$ws.Cells.Item(13, 13).Value = 0.6607217687215012
$ws.Cells.Item(13, 14).Value = 2.837474697957219
$ws.Cells.Item(14, 2).Value = 3.78463507652009
$ws.Cells.Item(14, 3).Value = 0.9624353811694277
$ws.Cells.Item(14, 4).Value = 0.1105975714622929
$ws.Cells.Item(14, 5).Value = 0.0357251329179924
$ws.Cells.Item(14, 7).Value = 0.002563323136809103
$ws.Cells.Item(14, 9).Value = 2.61152063895014
$ws.Cells.Item(14, 12).Value = 0.2512014200515011
$ws.Cells.Item(14, 13).Value = 0.6558744180170493
$ws.Cells.Item(14, 14).Value = 2.835730517465009
$ws.Cells.Item(15, 2).Value = 3.76438317476925
$ws.Cells.Item(15, 3).Value = 0.9540950297712811
$ws.Cells.Item(15, 4).Value = 0.1096222820417552
$ws.Cells.Item(15, 5).Value = 0.0354421406334815
$ws.Cells.Item(15, 7).Value = 0.00256412951405884
$ws.Cells.Item(15, 9).Value = 2.606085110352453
$ws.Cells.Item(15, 12).Value = 0.2506979581745838
$ws.Cells.Item(15, 13).Value = 0.6529133513601408
$ws.Cells.Item(15, 14).Value = 2.83467690637633
$ws.Cells.Item(16, 2).Value = 3.649050309812651
$ws.Cells.Item(16, 3).Value = 0.9064635076065883
$ws.Cells.Item(16, 4).Value = 0.1040516369648969
$ws.Cells.Item(16, 5).Value = 0.03381891263824954
$ws.Cells.Item(16, 7).Value = 0.002568818198651777
$ws.Cells.Item(16, 9).Value = 2.575323427547147
$ws.Cells.Item(16, 12).Value = 0.2478523472617695
$ws.Cells.Item(16, 13).Value = 0.636067798654409
$ws.Cells.Item(16, 14).Value = 2.828864060647902
$ws.Cells.Item(17, 2).Value = 3.578937242531083
$ws.Cells.Item(17, 3).Value = 0.8773867173488838
$ws.Cells.Item(17, 4).Value = 0.1006503212662437
$ws.Cells.Item(17, 5).Value = 0.03282166733990621
$ws.Cells.Item(17, 7).Value = 0.002571755113854874
$ws.Cells.Item(17, 9).Value = 2.556795476778575
$ws.Cells.Item(17, 12).Value = 0.2461418126684833
$ws.Cells.Item(17, 13).Value = 0.6258427586045059
$ws.Cells.Item(17, 14).Value = 2.825497626128083
$ws.Cells.Item(18, 2).Value = 3.538843310748859
$ws.Cells.Item(18, 3).Value = 0.8607143350986917
$ws.Cells.Item(18, 4).Value = 0.0986997812230328
$ws.Cells.Item(18, 5).Value = 0.0322475081280551
$ws.Cells.Item(18, 7).Value = 0.002573466677164761
$ws.Cells.Item(18, 9).Value = 2.546264041484591
$ws.Cells.Item(18, 12).Value = 0.2451708333703948
$ws.Cells.Item(18, 13).Value = 0.6200014258956301
$ws.Cells.Item(18, 14).Value = 2.823634210339435
$ws.Cells.Item(19, 2).Value = 3.525308075379485
$ws.Cells.Item(19, 3).Value = 0.8550781628232471
$ws.Cells.Item(19, 4).Value = 0.09804034904165349
$ws.Cells.Item(19, 5).Value = 0.0320530061782307
$ws.Cells.Item(19, 7).Value = 0.002574050024843093
$ws.Cells.Item(19, 9).Value = 2.542719729990509
$ws.Cells.Item(19, 12).Value = 0.2448442835728173
$ws.Cells.Item(19, 13).Value = 0.6180304675853563
$ws.Cells.Item(19, 14).Value = 2.82301576077424
$ws.Cells.Item(20, 2).Value = 3.58637671487071
$ws.Cells.Item(20, 3).Value = 0.8804766068572007
$ws.Cells.Item(20, 4).Value = 0.1010117937279063
$ws.Cells.Item(20, 5).Value = 0.03292788395093282
$ws.Cells.Item(20, 7).Value = 0.002571440165041494
$ws.Cells.Item(20, 9).Value = 2.55875481679989
$ws.Cells.Item(20, 12).Value = 0.24632256849047
$ws.Cells.Item(20, 13).Value = 0.6269271021101943
$ws.Cells.Item(20, 14).Value = 2.825848436722623
$ws.Cells.Item(21, 2).Value = 3.794359265113599
$ws.Cells.Item(21, 3).Value = 0.9664377357677267
$ws.Cells.Item(21, 4).Value = 0.1110655776253964
$ws.Cells.Item(21, 5).Value = 0.03586081000502972
$ws.Cells.Item(21, 7).Value = 0.002562937651420238
$ws.Cells.Item(21, 9).Value = 2.614134000504379
$ws.Cells.Item(21, 12).Value = 0.2514435444999776
$ws.Cells.Item(21, 13).Value = 0.657296517894622
$ws.Cells.Item(21, 14).Value = 2.836239749051714
$ws.Cells.Item(22, 2).Value = 3.932405999690673
$ws.Cells.Item(22, 3).Value = 1.023100205432343
$ws.Cells.Item(22, 4).Value = 0.1176903249757544
$ws.Cells.Item(22, 5).Value = 0.03777338567346433
$ws.Cells.Item(22, 7).Value = 0.002557578056365238
$ws.Cells.Item(22, 9).Value = 2.651461435747152
$ws.Cells.Item(22, 12).Value = 0.2549060067203186
$ws.Cells.Item(22, 13).Value = 0.67750538407752
$ws.Cells.Item(22, 14).Value = 2.843689974315424
$ws.Cells.Item(23, 2).Value = 3.85853203332158
$ws.Cells.Item(23, 3).Value = 0.9928132815356889
$ws.Cells.Item(23, 4).Value = 0.1141495167582889
$ws.Cells.Item(23, 5).Value = 0.03675295801242129
$ws.Cells.Item(23, 7).Value = 0.002560420595378643
$ws.Cells.Item(23, 9).Value = 2.63143445767362
$ws.Cells.Item(23, 12).Value = 0.2530473999686933
$ws.Cells.Item(23, 13).Value = 0.666686236393744
$ws.Cells.Item(23, 14).Value = 2.839652869112967
$ws.Cells.Item(24, 2).Value = 3.583012660184693
$ws.Cells.Item(24, 3).Value = 0.8790795302008974
$ws.Cells.Item(24, 4).Value = 0.1008483567234038
$ws.Cells.Item(24, 5).Value = 0.03287986601412385
$ws.Cells.Item(24, 7).Value = 0.002571582481349467
$ws.Cells.Item(24, 9).Value = 2.557868623945538
$ws.Cells.Item(24, 12).Value = 0.2462408100777935
$ws.Cells.Item(24, 13).Value = 0.6264367546719143
$ws.Cells.Item(24, 14).Value = 2.825689611123465
$ws.Cells.Item(25, 2).Value = 3.29471382866609
$ws.Cells.Item(25, 3).Value = 0.7583683289681744
$ws.Cells.Item(25, 4).Value = 0.08672154255728515
$ws.Cells.Item(25, 5).Value = 0.02867983450942901
$ws.Cells.Item(25, 7).Value = 0.002584476942991809
$ws.Cells.Item(25, 9).Value = 2.483298058301528
$ws.Cells.Item(25, 12).Value = 0.2393904739806914
$ws.Cells.Item(25, 13).Value = 0.584540837694199
$ws.Cells.Item(25, 14).Value = 2.813408059575139
